$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated cryptocurrency list (price / volume(1h) columns, and the
# Hedera/MXToken row order swap) as refreshed by the scheduled GitHub
# Actions job.
#
# Values are written through a small helper so that numeric-looking
# text such as "6.560" or "3.000" is preserved exactly as text (the
# source workbook stores these as inline strings), instead of being
# silently normalized into a number by Excel. The original cell style
# is captured and restored afterwards so no stray formatting/style is
# introduced by the temporary text NumberFormat.
function Set-TextValue($range, [string]$value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "27.525.80"
Set-TextValue $ws.Range("E2") "  +2.21%  "
Set-TextValue $ws.Range("D3") "1.868.63"
Set-TextValue $ws.Range("E3") "  +1.28%  "
Set-TextValue $ws.Range("E4") "  +0.23%  "
Set-TextValue $ws.Range("D5") "312.32"
Set-TextValue $ws.Range("E5") "  +0.87%  "
Set-TextValue $ws.Range("D6") "1.013"
Set-TextValue $ws.Range("E6") "  +0.30%  "
Set-TextValue $ws.Range("E7") "  +0.31%  "
Set-TextValue $ws.Range("D8") "0.3775"
Set-TextValue $ws.Range("E8") "  +2.94%  "
Set-TextValue $ws.Range("D9") "0.07365"
Set-TextValue $ws.Range("E9") "  +1.94%  "
Set-TextValue $ws.Range("D10") "0.9366"
Set-TextValue $ws.Range("E10") "  +1.16%  "
Set-TextValue $ws.Range("D11") "20.71"
Set-TextValue $ws.Range("E11") "  +5.35%  "
Set-TextValue $ws.Range("D12") "0.07842"
Set-TextValue $ws.Range("E12") "  +1.86%  "
Set-TextValue $ws.Range("D13") "1.882.70"
Set-TextValue $ws.Range("E13") "  +1.51%  "
Set-TextValue $ws.Range("D14") "5.454"
Set-TextValue $ws.Range("E14") "  +2.57%  "
Set-TextValue $ws.Range("D15") "6.560"
Set-TextValue $ws.Range("E15") "  +2.42%  "
Set-TextValue $ws.Range("D16") "90.67"
Set-TextValue $ws.Range("E16") "  +2.09%  "
Set-TextValue $ws.Range("E17") "  +0.20%  "
Set-TextValue $ws.Range("D18") "0.000008918"
Set-TextValue $ws.Range("E18") "  +3.21%  "
Set-TextValue $ws.Range("E19") "  +0.22%  "
Set-TextValue $ws.Range("D20") "14.90"
Set-TextValue $ws.Range("E20") "  +2.57%  "
Set-TextValue $ws.Range("D21") "27.533.85"
Set-TextValue $ws.Range("E21") "  +2.09%  "
Set-TextValue $ws.Range("D22") "5.128"
Set-TextValue $ws.Range("E22") "  +1.38%  "
Set-TextValue $ws.Range("D23") "10.71"
Set-TextValue $ws.Range("E23") "  +0.58%  "
Set-TextValue $ws.Range("D24") "1.953"
Set-TextValue $ws.Range("E24") "  +1.39%  "
Set-TextValue $ws.Range("D25") "154.20"
Set-TextValue $ws.Range("E25") "  +1.18%  "
Set-TextValue $ws.Range("D26") "18.49"
Set-TextValue $ws.Range("E26") "  +1.88%  "
Set-TextValue $ws.Range("D27") "2.025"
Set-TextValue $ws.Range("E27") "  +1.64%  "
Set-TextValue $ws.Range("D28") "115.85"
Set-TextValue $ws.Range("E28") "  +1.51%  "
Set-TextValue $ws.Range("D29") "4.993"
Set-TextValue $ws.Range("E29") "  +1.09%  "
Set-TextValue $ws.Range("D30") "0.08920"
Set-TextValue $ws.Range("E30") "  +0.50%  "
Set-TextValue $ws.Range("D31") "3.333"
Set-TextValue $ws.Range("E31") "  +0.44%  "
Set-TextValue $ws.Range("D32") "1.218"
Set-TextValue $ws.Range("E32") "  +4.04%  "
Set-TextValue $ws.Range("D33") "0.7602"
Set-TextValue $ws.Range("E33") "  +2.25%  "
Set-TextValue $ws.Range("D34") "4.619"
Set-TextValue $ws.Range("E34") "  +2.89%  "
Set-TextValue $ws.Range("E35") "  +0.29%  "
Set-TextValue $ws.Range("D36") "0.02056"
Set-TextValue $ws.Range("E36") "  +5.02%  "
Set-TextValue $ws.Range("E37") "  -0.16%  "
Set-TextValue $ws.Range("B38") "MXToken"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-TextValue $ws.Range("D38") "3.000"
Set-TextValue $ws.Range("E38") "  +0.66%  "
Set-TextValue $ws.Range("B39") "Hedera"
Set-TextValue $ws.Range("C39") "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-TextValue $ws.Range("D39") "0.05274"
Set-TextValue $ws.Range("E39") "  +0.17%  "
Set-TextValue $ws.Range("D40") "0.5336"
Set-TextValue $ws.Range("E40") "  +2.92%  "
Set-TextValue $ws.Range("D41") "7.085"
Set-TextValue $ws.Range("E41") "  +1.42%  "
Set-TextValue $ws.Range("D42") "8.501"
Set-TextValue $ws.Range("E42") "  +3.85%  "
Set-TextValue $ws.Range("D43") "0.1525"
Set-TextValue $ws.Range("E43") "  +1.05%  "
Set-TextValue $ws.Range("D44") "10.69"
Set-TextValue $ws.Range("E44") "  +1.04%  "
Set-TextValue $ws.Range("D45") "0.4807"
Set-TextValue $ws.Range("E45") "  +2.00%  "
Set-TextValue $ws.Range("E46") "  +0.27%  "
Set-TextValue $ws.Range("D47") "1.660"
Set-TextValue $ws.Range("E47") "  +3.70%  "
Set-TextValue $ws.Range("D48") "102.96"
Set-TextValue $ws.Range("E48") "  +1.55%  "
Set-TextValue $ws.Range("D49") "67.48"
Set-TextValue $ws.Range("E49") "  +2.58%  "
Set-TextValue $ws.Range("D50") "0.06083"
Set-TextValue $ws.Range("E50") "  +1.04%  "
Set-TextValue $ws.Range("D51") "0.9263"
Set-TextValue $ws.Range("E51") "  +4.62%  "
